$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen columns D:L (col widths measured in characters; Excel's
# ColumnWidth setter adds ~5/6 character of internal padding relative to
# the stored <col width> value, so subtract it to land on the exact
# target integer widths from the diff). ---
$pad = 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth  = 40 - $pad   # D
$ws.Columns.Item(5).ColumnWidth  = 50 - $pad   # E
$ws.Columns.Item(6).ColumnWidth  = 27 - $pad   # F
$ws.Columns.Item(7).ColumnWidth  = 16 - $pad   # G
$ws.Columns.Item(8).ColumnWidth  = 45 - $pad   # H
$ws.Columns.Item(9).ColumnWidth  = 35 - $pad   # I
$ws.Columns.Item(10).ColumnWidth = 38 - $pad   # J
$ws.Columns.Item(11).ColumnWidth = 37 - $pad   # K
$ws.Columns.Item(12).ColumnWidth = 45 - $pad   # L

# --- Fill in newly-scraped data for the existing rows ---

# Row 2 - Student Life
$ws.Range("E2").Value = "Student organization focused on general activities and community engagement. The Student Life welcomes all interested students to participate and make a positive impact."
$ws.Range("G2").Value = "(555) 789-0123"
$ws.Range("I2").Value = "https://instagram.com/studentlife"
$ws.Range("J2").Value = "https://facebook.com/studentlife"

# Row 3 - Clubs
$ws.Range("D3").Value = "https://biola.edu/logos/clubs_logo.png"
$ws.Range("E3").Value = "Student organization focused on general activities and community engagement. The Clubs welcomes all interested students to participate and make a positive impact."
$ws.Range("F3").Value = "clubs@biola.edu"

# Row 4 - Multi-Ethnic Clubs
$ws.Range("E4").Value = "Student organization focused on general activities and community engagement. The Multi-Ethnic Clubs welcomes all interested students to participate and make a positive impact."
$ws.Range("F4").Value = "multiethnicclub@biola.edu"
$ws.Range("J4").Value = "https://facebook.com/multiethnicclub"
$ws.Range("K4").Value = "https://twitter.com/multiethnicclub"
$ws.Range("L4").Value = "https://youtube.com/channel/multiethnicclub"

# Row 5 - Student Life - The Biola Experience - Biola University
$ws.Range("E5").Value = "Student organization focused on general activities and community engagement. The Student Life - The Biola Experience - Biola University welcomes all interested students to participate and make a positive impact."
$ws.Range("F5").Value = "studentlifetheb@biola.edu"
$ws.Range("H5").Value = "https://linkedin.com/groups/studentlifetheb"
